$d = $word.ActiveDocument

# Locate the two target list-paragraphs by their text content so the script
# is resilient to any incidental paragraph-index shifts.
$translateIdx = -1
$prepareIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*translate stimuli to German*") {
        $translateIdx = $i
    }
    if ($t -like "*Prepare an initial proposal for an experiment*") {
        $prepareIdx = $i
    }
}

# --- Change 1: split the "translate stimuli" paragraph into multiple runs,
# wrapping the German words with spell-check proofErr boundary markers
# (reisen / anreisen / nach / nach), exactly as Word's proofing pass does
# for words it does not recognise. ---
$p1 = $d.Paragraphs($translateIdx)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">To do: translate stimuli to German and ask Michael about (1) </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>reisen</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> or </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>anreisen</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> and (2) </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>nach</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> or </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>nach</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Afrika</w:t></w:r>' + `
  '</w:p>'
$p1.Range.InsertXML($xml1) | Out-Null

# --- Change 2: merge the "To do: " / "Prepare an initial proposal..." runs
# into a single run. A same-text Find/Replace that spans both runs coalesces
# them into one run without altering the visible text. ---
$d.Content.Find.Execute("To do: Prepare an initial proposal for an experiment", $true, $false, $false, $false, $false, $true, 1, $false, "To do: Prepare an initial proposal for an experiment", 2) | Out-Null

# --- Change 3: add a new to-do bullet right after the "Prepare an initial
# proposal" item, using the same list style/numbering. ---
$p2 = $d.Paragraphs($prepareIdx)
$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = "To do: discuss number and style of stimuli with Michael/James"
